$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'68.841.32"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = "'3.926.12"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.69%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'606.15"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = "'165.95"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('D7').Value = "'3.923.19"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.72%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').Value = "'37.33"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = "'4.581.45"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.62%  '
$ws.Range('D16').Value = "'3.882.81"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.07%  '
$ws.Range('D17').Value = "'68.981.76"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = "'7.49"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = "'17.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.27%  '
$ws.Range('D21').Value = "'11.14"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').Value = "'488.31"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = "'0.724"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +12.05%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').Value = "'2.27"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = "'12.13"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('D28').Value = "'10.17"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.23%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').Value = "'4.076.53"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.47%  '
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').Value = "'7.89"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.41%  '
$ws.Range('D34').Value = "'32.33"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('D35').Value = "'3.873.68"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.05%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  +3.31%  '
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('D39').Value = "'5.93"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = "'1.00"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = "'0.322"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('D42').Value = "'440.16"
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = "'3.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.20%  '
$ws.Range('D44').Value = "'2.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').Value = "'48.45"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = "'8.54"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = "'2.851.30"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('D49').Value = "'26.25"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.99%  '
$ws.Range('D50').Value = "'141.77"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = "'0.0357"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.35%  '
